$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("Data")

# --- Update text on the "About" sheet ---
$wsAbout.Range("B6").Value = "January 2020 and November 2020"
$wsAbout.Range("A27").Value = "As of EPS 3.1, this variable is set up to model the impacts of the 2020"
$wsAbout.Range("A28").Value = "SARS-CoV-2 pandemic.  It uses the latest data available as of November 10,"

# --- Update text + values on the "Data" sheet ---
$wsData.Range("A3").Value = "November STEO"

$wsData.Range("B3").Value = 19092
$wsData.Range("C3").Value = 18411
$wsData.Range("D3").Value = 19098

# --- Update selection on the Data sheet ---
$wsData.Range("B12").Select()

# Restore the originally active sheet/tab
$wsAbout.Activate()

$wb.Save()
